$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.729.88"
$ws.Range("E2").Value = "  +5.27%  "
$ws.Range("D3").Value = "2.300.63"
$ws.Range("E3").Value = "  +3.56%  "
$ws.Range("E4").Value = "  -0.65%  "
$ws.Range("D5").Value = "'302.38"
$ws.Range("E5").Value = "  +1.30%  "
$ws.Range("D6").Value = "'101.53"
$ws.Range("E6").Value = "  +12.74%  "
$ws.Range("D7").Value = "'0.569"
$ws.Range("E7").Value = "  +0.97%  "
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("D9").Value = "'0.531"
$ws.Range("E9").Value = "  +8.41%  "
$ws.Range("D10").Value = "'36.81"
$ws.Range("E10").Value = "  +11.30%  "
$ws.Range("E11").Value = "  +2.06%  "
$ws.Range("D12").Value = "'7.36"
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("D14").Value = "2.653.63"
$ws.Range("E14").Value = "  +3.58%  "
$ws.Range("D15").Value = "2.302.78"
$ws.Range("E15").Value = "  +3.62%  "
$ws.Range("D16").Value = "'14.01"
$ws.Range("E16").Value = "  +3.16%  "
$ws.Range("D17").Value = "'0.820"
$ws.Range("D18").Value = "46.693.82"
$ws.Range("E18").Value = "  +5.83%  "
$ws.Range("D19").Value = "'13.45"
$ws.Range("E19").Value = "  +21.48%  "
$ws.Range("E20").Value = "  +4.25%  "
$ws.Range("E21").Value = "  +3.13%  "
$ws.Range("E22").Value = "  +3.77%  "
$ws.Range("D23").Value = "'248.44"
$ws.Range("E23").Value = "  +5.43%  "
$ws.Range("E24").Value = "  +5.77%  "
$ws.Range("E25").Value = "  +5.93%  "
$ws.Range("E26").Value = "  -1.20%  "
$ws.Range("D27").Value = "'44.89"
$ws.Range("E27").Value = "  +16.59%  "
$ws.Range("E28").Value = "  +1.23%  "
$ws.Range("E29").Value = "  +5.77%  "
$ws.Range("D30").Value = "'20.13"
$ws.Range("E30").Value = "  +3.14%  "
$ws.Range("E31").Value = "  +7.81%  "
$ws.Range("D32").Value = "'147.66"
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("E33").Value = "  +6.60%  "
$ws.Range("E34").Value = "  +3.41%  "
$ws.Range("D35").Value = "'3.15"
$ws.Range("E35").Value = "  +11.16%  "
$ws.Range("E36").Value = "  +9.95%  "
$ws.Range("E37").Value = "  +2.66%  "
$ws.Range("D38").Value = "'1.80"
$ws.Range("E38").Value = "  +7.94%  "
$ws.Range("D39").Value = "'15.90"
$ws.Range("D40").Value = "'4.04"
$ws.Range("E40").Value = "  +13.37%  "
$ws.Range("D41").Value = "'3.51"
$ws.Range("E41").Value = "  +10.07%  "
$ws.Range("E42").Value = "  +0.86%  "
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("D44").Value = "1.868.31"
$ws.Range("E44").Value = "  +2.32%  "
$ws.Range("E45").Value = "  +11.16%  "
$ws.Range("D46").Value = "'87.85"
$ws.Range("E46").Value = "  +19.11%  "
$ws.Range("E47").Value = "  +9.87%  "
$ws.Range("D48").Value = "'74.25"
$ws.Range("E48").Value = "  +11.01%  "
$ws.Range("D49").Value = "'4.90"
$ws.Range("E49").Value = "  +8.53%  "
$ws.Range("D50").Value = "'97.12"
$ws.Range("E50").Value = "  +3.20%  "
$ws.Range("D51").Value = "'8.07"
$ws.Range("E51").Value = "  +5.82%  "
